$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (Coin name) swaps for rows 38/39 ---
$ws.Range("B38").Value = "VeChain"
$ws.Range("B39").Value = "LidoDAOToken"

# --- Column C (Link) swaps for rows 38/39 ---
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"

# --- Column D values that remain non-numeric text (contain multiple dots) ---
# These can be assigned directly as strings without any special handling.
$ws.Range("D2").Value = "42.095.23"
$ws.Range("D3").Value = "2.294.07"
$ws.Range("D16").Value = "2.644.63"
$ws.Range("D17").Value = "2.318.14"
$ws.Range("D18").Value = "42.050.02"

# --- Column D values that look like plain numbers ---
# Excel would normally auto-convert these to a Number type. The source
# workbook stores Price as text (inline string) for every row, so we force
# text entry with a leading apostrophe, then strip the resulting "quote
# prefix" / text-number-format by pasting the (always-plain) format from
# D4, a cell whose Price value is never modified by this update.
$ws.Range("D5").Value = "'314.22"
$ws.Range("D6").Value = "'105.79"
$ws.Range("D7").Value = "'0.626"
$ws.Range("D9").Value = "'0.608"
$ws.Range("D10").Value = "'39.86"
$ws.Range("D11").Value = "'0.0911"
$ws.Range("D12").Value = "'8.36"
$ws.Range("D14").Value = "'0.973"
$ws.Range("D15").Value = "'15.35"
$ws.Range("D19").Value = "'7.63"
$ws.Range("D21").Value = "'72.71"
$ws.Range("D22").Value = "'3.52"
$ws.Range("D23").Value = "'258.40"
$ws.Range("D25").Value = "'9.75"
$ws.Range("D27").Value = "'10.96"
$ws.Range("D29").Value = "'22.71"
$ws.Range("D30").Value = "'36.24"
$ws.Range("D31").Value = "'165.59"
$ws.Range("D32").Value = "'0.0888"
$ws.Range("D35").Value = "'0.119"
$ws.Range("D37").Value = "'4.62"
$ws.Range("D38").Value = "'0.0351"
$ws.Range("D39").Value = "'2.90"
$ws.Range("D40").Value = "'3.61"
$ws.Range("D41").Value = "'98.46"
$ws.Range("D42").Value = "'1.48"
$ws.Range("D43").Value = "'70.75"
$ws.Range("D44").Value = "'0.227"
$ws.Range("D46").Value = "'12.18"
$ws.Range("D47").Value = "'112.56"
$ws.Range("D48").Value = "'78.19"
$ws.Range("D49").Value = "'9.13"
$ws.Range("D50").Value = "'5.32"
$ws.Range("D4").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("D21").PasteSpecial(-4122)
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("D32").PasteSpecial(-4122)
$ws.Range("D35").PasteSpecial(-4122)
$ws.Range("D37").PasteSpecial(-4122)
$ws.Range("D38").PasteSpecial(-4122)
$ws.Range("D39").PasteSpecial(-4122)
$ws.Range("D40").PasteSpecial(-4122)
$ws.Range("D41").PasteSpecial(-4122)
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("D43").PasteSpecial(-4122)
$ws.Range("D44").PasteSpecial(-4122)
$ws.Range("D46").PasteSpecial(-4122)
$ws.Range("D47").PasteSpecial(-4122)
$ws.Range("D48").PasteSpecial(-4122)
$ws.Range("D49").PasteSpecial(-4122)
$ws.Range("D50").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Column E (Volume(1h)) percentage updates ---
$ws.Range("E2").Value = "  -1.46%  "
$ws.Range("E3").Value = "  -1.99%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("E6").Value = "  +1.02%  "
$ws.Range("E7").Value = "  -1.79%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -1.18%  "
$ws.Range("E10").Value = "  -1.65%  "
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("E15").Value = "  -3.18%  "
$ws.Range("E16").Value = "  -1.80%  "
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("E18").Value = "  -1.40%  "
$ws.Range("E19").Value = "  -1.44%  "
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("E21").Value = "  -5.91%  "
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("E25").Value = "  +1.10%  "
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("E27").Value = "  -3.61%  "
$ws.Range("E28").Value = "  +2.56%  "
$ws.Range("E29").Value = "  -1.49%  "
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("E31").Value = "  -5.33%  "
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("E33").Value = "  -2.15%  "
$ws.Range("E34").Value = "  -2.92%  "
$ws.Range("E35").Value = "  +6.46%  "
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("E37").Value = "  +2.01%  "
$ws.Range("E38").Value = "  -1.03%  "
$ws.Range("E39").Value = "  +10.38%  "
$ws.Range("E40").Value = "  -3.44%  "
$ws.Range("E41").Value = "  +15.48%  "
$ws.Range("E42").Value = "  +0.88%  "
$ws.Range("E43").Value = "  +0.51%  "
$ws.Range("E44").Value = "  -2.24%  "
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("E46").Value = "  +3.20%  "
$ws.Range("E47").Value = "  -2.29%  "
$ws.Range("E48").Value = "  +6.96%  "
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("E50").Value = "  -3.35%  "
$ws.Range("E51").Value = "  +2.41%  "

